$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 ("Save"), matching the style of the other header cells
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add new "Save" data column values for rows 2-6
$values = @(0, 1, 1, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
